# teacher_init.xlsx update: return/validate file import teacher
# - Fill in the missing "maso" (B7) and "hodem" (C7) values for the 6th
#   teacher row (Nguyen Thi Thanh Binh), which were left blank before.
# - Re-apply a Vietnamese-locale short date format (d/m/yyyy) to the
#   "ngay_sinh" column instead of the default m/d/yyyy format.
# - Remove the stray trailing blank row left at the bottom of the sheet.
# - Restore the selected cell / window state left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fill in the previously blank cells on row 7 ---------------------
$ws.Cells.Item(7, 2).Value = 20012016
$ws.Cells.Item(7, 3).Value = "Nguyễn Thị"

# --- apply the Vietnamese date format to the whole ngay_sinh column --
$ws.Columns.Item(6).NumberFormat = "[$-1010000]d/m/yyyy;@"

# --- the Kim Chi e-mail cell reverted to the plain default style ------
$ws.Range("E3").Style = "Normal"

# --- drop the stray empty row 12 at the bottom of the table -----------
$ws.Rows.Item(12).Delete()

# --- restore the view/selection state ---------------------------------
[void]$ws.Range("G19").Select()
